$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.753.01"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").Value = "2.107.35"
$ws.Range("E3").Value = "  +2.21%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'233.95"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("E6").Value = "  +0.85%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'57.71"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +1.62%  "

$ws.Range("E10").Value = "  +2.40%  "

$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("D12").Value = "2.418.49"
$ws.Range("E12").Value = "  +2.26%  "

$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("E14").Value = "  +2.06%  "

$ws.Range("E15").Value = "  +0.93%  "

$ws.Range("D16").Value = "'5.20"
$ws.Range("E16").Value = "  +0.99%  "

$ws.Range("D17").Value = "2.141.23"
$ws.Range("E17").Value = "  +3.92%  "

$ws.Range("D18").Value = "37.680.86"
$ws.Range("E18").Value = "  +1.29%  "

$ws.Range("D19").Value = "'6.21"
$ws.Range("E19").Value = "  -2.84%  "

$ws.Range("D20").Value = "'70.25"
$ws.Range("E20").Value = "  +1.72%  "

$ws.Range("E21").Value = "  +1.22%  "

$ws.Range("D22").Value = "'227.06"
$ws.Range("E22").Value = "  +1.07%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "'2.41"
$ws.Range("E24").Value = "  -0.94%  "

$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  +1.07%  "

$ws.Range("D26").Value = "'169.48"
$ws.Range("E26").Value = "  +2.40%  "

$ws.Range("D27").Value = "'8.96"
$ws.Range("E27").Value = "  +2.00%  "

$ws.Range("D28").Value = "'0.133"
$ws.Range("E28").Value = "  +5.09%  "

$ws.Range("E29").Value = "  -2.04%  "

$ws.Range("D30").Value = "'19.42"
$ws.Range("E30").Value = "  +1.93%  "

$ws.Range("E31").Value = "  +1.15%  "

$ws.Range("D32").Value = "'4.62"
$ws.Range("E32").Value = "  +2.90%  "

$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").Value = "'2.58"
$ws.Range("E33").Value = "  +1.52%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0622"
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("E35").Value = "  +0.30%  "

$ws.Range("D36").Value = "'3.43"
$ws.Range("E36").Value = "  +5.02%  "

$ws.Range("E37").Value = "  +4.10%  "

$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("D39").Value = "'5.42"
$ws.Range("E39").Value = "  -6.98%  "

$ws.Range("D40").Value = "'0.101"
$ws.Range("E40").Value = "  +8.61%  "

$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").Value = "'96.54"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").Value = "1.458.93"
$ws.Range("E43").Value = "  -2.13%  "

$ws.Range("E44").Value = "  +1.30%  "

$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("D46").Value = "'4.09"
$ws.Range("E46").Value = "  -11.97%  "

$ws.Range("E47").Value = "  +2.90%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'15.39"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'3.03"
$ws.Range("E49").Value = "  +2.53%  "

$ws.Range("E50").Value = "  +1.81%  "

$ws.Range("D51").Value = "2.303.90"
$ws.Range("E51").Value = "  +2.32%  "
